$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumen")

$ws.Range("H3").Value = 279
$ws.Range("I3").Value = 71
$ws.Range("J3").Value = 3753
$ws.Range("K3").Value = 111.4
$ws.Range("L3").Value = 44
$ws.Range("M3").Value = 21
$ws.Range("N3").Value = 6
$ws.Range("Q3").Value = 97

$ws.Range("H4").Value = 132
$ws.Range("I4").Value = 32
$ws.Range("J4").Value = 679
$ws.Range("K4").Value = 43.4
$ws.Range("L4").Value = 25
$ws.Range("M4").Value = 7
$ws.Range("Q4").Value = 51

$ws.Range("H5").Value = 179
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 1419
$ws.Range("K5").Value = 75.40000000000001
$ws.Range("L5").Value = 44
$ws.Range("M5").Value = 18
$ws.Range("N5").Value = 8
$ws.Range("Q5").Value = 78

$ws.Range("H6").Value = 80
$ws.Range("I6").Value = 19
$ws.Range("J6").Value = 761
$ws.Range("K6").Value = 58
$ws.Range("L6").Value = 11
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 3
$ws.Range("Q6").Value = 28

$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 14
$ws.Range("J7").Value = 144
$ws.Range("K7").Value = 59.8
$ws.Range("L7").Value = 7
$ws.Range("M7").Value = 7
$ws.Range("N7").Value = 0
$ws.Range("Q7").Value = 18

$ws.Range("H8").Value = 128
$ws.Range("I8").Value = 20
$ws.Range("J8").Value = 646
$ws.Range("K8").Value = 63.3
$ws.Range("L8").Value = 8
$ws.Range("M8").Value = 9
$ws.Range("N8").Value = 3
$ws.Range("Q8").Value = 36

$ws.Range("H9").Value = 169
$ws.Range("I9").Value = 61
$ws.Range("J9").Value = 3654
$ws.Range("K9").Value = 67.40000000000001
$ws.Range("L9").Value = 38
$ws.Range("M9").Value = 19
$ws.Range("N9").Value = 4
$ws.Range("Q9").Value = 67

$ws.Range("H10").Value = 41
$ws.Range("I10").Value = 13
$ws.Range("J10").Value = 489
$ws.Range("K10").Value = 72
$ws.Range("L10").Value = 12
$ws.Range("M10").Value = 1
$ws.Range("Q10").Value = 14

$ws.Range("H11").Value = 39
$ws.Range("I11").Value = 14
$ws.Range("J11").Value = 921
$ws.Range("K11").Value = 46.5
$ws.Range("L11").Value = 12
$ws.Range("M11").Value = 2
$ws.Range("N11").Value = 0
$ws.Range("Q11").Value = 15

$ws.Range("H12").Value = 34
$ws.Range("I12").Value = 11
$ws.Range("J12").Value = 236
$ws.Range("K12").Value = 62.3
$ws.Range("L12").Value = 10
$ws.Range("M12").Value = 1
$ws.Range("Q12").Value = 11

$ws.Range("H13").Value = 34
$ws.Range("I13").Value = 7
$ws.Range("J13").Value = 448
$ws.Range("K13").Value = 26.9
$ws.Range("N13").Value = 7
$ws.Range("Q13").Value = 7
